$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers in row 2 (entry order controls the shared-string table order) ---
$ws.Range("E2").Value = "Min"
$ws.Range("G2").Value = "Count"
$ws.Range("F2").Value = "Max"

# New headers pick up the same bold style already used for B2:D2
$ws.Range("E2:G2").Font.Bold = $true

# --- Row 3: new summary formulas (StdDev/Sum in C3:D3 already existed) ---
$ws.Range("E3").Formula = "=MIN(B:B)"
$ws.Range("F3").Formula = "=MAX(B:B)"
$ws.Range("G3").Formula = "=COUNT(B:B)"

# Min/Max share the scientific-notation format already applied to StdDev/Sum
$ws.Range("E3").NumberFormat = "0.00000000000000000E+00"
$ws.Range("F3").NumberFormat = "0.00000000000000000E+00"

# --- Row 4: mirrors of the row-3 summary cells, plus the first tracked value ---
$ws.Range("B4").Value = 0.1
$ws.Range("E4").Formula = "=E3"
$ws.Range("F4").Formula = "=F3"

# Mirrors share the long-decimal format already applied to the StdDev/Sum mirrors
$ws.Range("E4").NumberFormat = "0.00000000000000000"
$ws.Range("F4").NumberFormat = "0.00000000000000000"

# --- Remaining tracked values replace the old 0.17 / 0.32 / -0.15 / 1.07 set ---
$ws.Range("B5").Value = 0.2
$ws.Range("B6").Value = 0.3
$ws.Range("B7").Value = 0.4
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()

# --- Column widths: C:F widened to fit the new headers, G newly sized ---
$ws.Range("C:F").ColumnWidth = 34.5859375
$ws.Range("G:G").ColumnWidth = 20.1171875

# --- Selection moves from B7 to B4 ---
$ws.Range("B4").Select()
